# Generate Report for Handback
#
# The localization pipeline finished a handback pass: the 1d7ee5ab file
# (previously "Ready for handoff") is now "Handed back: in sync with en-US",
# same as 29f983bc. The report rows are re-sorted by status (handed-back
# rows first), so the 1d7ee5ab row moves from row 3 up to row 2 and the
# 29f983bc row moves from row 2 down to row 3 on every sheet. The
# "Latest Handback DateTime" for the zh-cn / de-de hand-back rows is
# refreshed to reflect the new handback event.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# "Overview" sheet: File Name | zh-cn | de-de
# ---------------------------------------------------------------------
$ov = $wb.Worksheets.Item("Overview")

$ov.Range("A2").Value = "1d7ee5ab-25dc-427c-9ba7-d760e14e9d1b.md"
$ov.Range("B2").Value = "Handed back: in sync with en-US"
$ov.Range("C2").Value = "Handed back: in sync with en-US"

$ov.Range("A3").Value = "29f983bc-ed00-4945-b663-5fc111d2269f.md"
$ov.Range("B3").Value = "Handed back: in sync with en-US"
$ov.Range("C3").Value = "Handed back: in sync with en-US"

# ---------------------------------------------------------------------
# "zh-cn" sheet
# ---------------------------------------------------------------------
$zh = $wb.Worksheets.Item("zh-cn")

$zh.Range("A2").Value = "1d7ee5ab-25dc-427c-9ba7-d760e14e9d1b.md"
$zh.Range("B2").Value = "Handed back: in sync with en-US"
$zh.Range("C2").Value = "1d7ee5ab-25dc-427c-9ba7-d760e14e9d1b.63194497c13f7696af8acce6728b295184dea86a.zh-cn.xlf"
$zh.Range("D2").Value = "2016-03-07 05:12:41"
$zh.Range("E2").Value = "1d7ee5ab-25dc-427c-9ba7-d760e14e9d1b.md"
$zh.Range("F2").Value = "1d7ee5ab-25dc-427c-9ba7-d760e14e9d1b.63194497c13f7696af8acce6728b295184dea86a.zh-cn.xlf"
$zh.Range("G2").Value = "2016-03-07 05:13:37"
$zh.Range("H2").Value = "Include"

$zh.Range("A3").Value = "29f983bc-ed00-4945-b663-5fc111d2269f.md"
$zh.Range("B3").Value = "Handed back: in sync with en-US"
$zh.Range("C3").Value = "29f983bc-ed00-4945-b663-5fc111d2269f.003e98a0eeedcc1d9f6f18cc32f67ddd27321de5.zh-cn.xlf"
$zh.Range("D3").Value = "2016-03-07 05:12:41"
$zh.Range("E3").Value = "29f983bc-ed00-4945-b663-5fc111d2269f.md"
$zh.Range("F3").Value = "29f983bc-ed00-4945-b663-5fc111d2269f.003e98a0eeedcc1d9f6f18cc32f67ddd27321de5.zh-cn.xlf"
$zh.Range("G3").Value = "2016-03-07 05:13:37"
$zh.Range("H3").Value = "Include"

# ---------------------------------------------------------------------
# "de-de" sheet
# ---------------------------------------------------------------------
$de = $wb.Worksheets.Item("de-de")

$de.Range("A2").Value = "1d7ee5ab-25dc-427c-9ba7-d760e14e9d1b.md"
$de.Range("B2").Value = "Handed back: in sync with en-US"
$de.Range("C2").Value = "1d7ee5ab-25dc-427c-9ba7-d760e14e9d1b.63194497c13f7696af8acce6728b295184dea86a.de-de.xlf"
$de.Range("D2").Value = "2016-03-07 05:12:52"
$de.Range("E2").Value = "1d7ee5ab-25dc-427c-9ba7-d760e14e9d1b.md"
$de.Range("F2").Value = "1d7ee5ab-25dc-427c-9ba7-d760e14e9d1b.63194497c13f7696af8acce6728b295184dea86a.de-de.xlf"
$de.Range("G2").Value = "2016-03-07 05:13:57"
$de.Range("H2").Value = "Include"

$de.Range("A3").Value = "29f983bc-ed00-4945-b663-5fc111d2269f.md"
$de.Range("B3").Value = "Handed back: in sync with en-US"
$de.Range("C3").Value = "29f983bc-ed00-4945-b663-5fc111d2269f.003e98a0eeedcc1d9f6f18cc32f67ddd27321de5.de-de.xlf"
$de.Range("D3").Value = "2016-03-07 05:12:52"
$de.Range("E3").Value = "29f983bc-ed00-4945-b663-5fc111d2269f.md"
$de.Range("F3").Value = "29f983bc-ed00-4945-b663-5fc111d2269f.003e98a0eeedcc1d9f6f18cc32f67ddd27321de5.de-de.xlf"
$de.Range("G3").Value = "2016-03-07 05:13:57"
$de.Range("H3").Value = "Include"
